$d = $word.ActiveDocument

# The "Tim Penguji" roles/placeholders are being rotated:
#   Penguji 1 (Ketua Tim Penguji) -> Ketua Sidang
#   Penguji 2 (Anggota)           -> Penguji 1
#   Penguji 3 (Pembimbing/Anggota)-> Penguji 2
# so every {Skor_...} / {Nama...} merge field tied to those roles is renamed
# to match (processed in an order that avoids one rename's output being
# re-matched by a later search).

$d.Content.Find.Execute("{Skor_Penguji1}", $true, $false, $false, $false, $false, $true, 1, $false, "{Skor_KetuaSidang}", 2)
$d.Content.Find.Execute("{Skor_Penguji2}", $true, $false, $false, $false, $false, $true, 1, $false, "{Skor_Penguji1}", 2)
$d.Content.Find.Execute("{Skor_Pembimbing}", $true, $false, $false, $false, $false, $true, 1, $false, "{Skor_Penguji2}", 2)

$d.Content.Find.Execute("{NamaPenguji1}", $true, $false, $false, $false, $false, $true, 1, $false, "{NamaKetuaSidang}", 2)
$d.Content.Find.Execute("{NamaPenguji2}", $true, $false, $false, $false, $false, $true, 1, $false, "{NamaPenguji1}", 2)
$d.Content.Find.Execute("{NamaPembimbing}", $true, $false, $false, $false, $false, $true, 1, $false, "{NamaPenguji2}", 2)
